$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the text of the shared string "-198.2% & -63.2%" -> "-198.2 & -63.2"
# (row 45, column L holds that result string)
$ws.Range("L45").Value = "-198.2 & -63.2"

# Row 45, column P: style only change (s=6 -> s=7), value unchanged.
# Copy P44's format (it already carries style 7) onto P45 without touching
# its value.
$ws.Range("P44").Copy()
$ws.Range("P45").PasteSpecial(-4122)

# Row 46: populate with a new run, mirroring the A-H/I values already used by
# rows 44-45, with new L/M/P results and update the row height to 19.5.
$ws.Range("A46").Value = "ukb51139_subset.csv"
$ws.Range("B46").Value = "28012 x 1081"
$ws.Range("C46").Value = "all"
$ws.Range("D46").Value = "no events"
$ws.Range("E46").Value = "> 140/80"
$ws.Range("F46").Value = "zscore"
$ws.Range("G46").Value = "median"
$ws.Range("H46").Value = "none"
$ws.Range("I46").Value = 25
$ws.Range("K46").Value = "N/A"
$ws.Range("L46").Value = "-55.6 & -27.9"
$ws.Range("M46").Value = "13.7 & 11.3"
$ws.Range("N46").Value = "N/A"
$ws.Range("O46").Value = "N/A"
$ws.Range("P46").Value = "20 batches"

$ws.Rows.Item(46).RowHeight = 19.5
$ws.Rows.Item(47).RowHeight = 19.5
$ws.Rows.Item(48).RowHeight = 19.5
